$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 4340
$ws.Range("I34").Value = 3012
$ws.Range("K34").Value = 3012
$ws.Range("M34").Value = -2809

$ws.Range("H36").Value = 4340
$ws.Range("I36").Value = 3012
$ws.Range("K36").Value = 3012
$ws.Range("M36").Value = -2297

$ws.Range("H100").Value = 333335330
$ws.Range("I100").Value = 500000480
$ws.Range("K100").Value = 500000480
$ws.Range("M100").Value = -499999939

$ws.Range("H112").Value = 2925002.8
$ws.Range("J112").Value = 2925002.8
$ws.Range("L112").Value = 8775008.399999999
$ws.Range("N112").Value = -8777224.399999999

$ws.Range("H129").Value = 162315.1
$ws.Range("J129").Value = 173488.72
$ws.Range("L129").Value = 520466.16
$ws.Range("N129").Value = -530466.16

$ws.Range("H138").Value = 2055.8235
$ws.Range("I138").Value = 1002
$ws.Range("K138").Value = 3006
$ws.Range("M138").Value = 2134

$ws.Range("H141").Value = 1032.1111
$ws.Range("I141").Value = 774.59186
$ws.Range("K141").Value = 2323.77558
$ws.Range("M141").Value = 2856.22442

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1550.238
$ws.Range("I2").Value = 1409.2059
$ws.Range("J2").Value = 2149.625
$ws.Range("K2").Value = 1409.2059
$ws.Range("L2").Value = 2149.625
$ws.Range("M2").Value = -1296.2059
$ws.Range("N2").Value = -2375.625

$ws.Range("H32").Value = 18020.166
$ws.Range("I32").Value = 18615.484
$ws.Range("K32").Value = 18615.484
$ws.Range("M32").Value = -18328.484

$ws.Range("H88").Value = 201110.8
$ws.Range("J88").Value = 1000054
$ws.Range("L88").Value = 1000054
$ws.Range("N88").Value = -1000866

$ws.Range("H91").Value = 201110.8
$ws.Range("J91").Value = 1000054
$ws.Range("L91").Value = 1000054
$ws.Range("N91").Value = -1002862

$ws.Range("H102").Value = 1600
$ws.Range("I102").Value = 1200
$ws.Range("K102").Value = 1200
$ws.Range("M102").Value = 422

$ws.Range("H116").Value = 1550.238
$ws.Range("I116").Value = 1409.2059
$ws.Range("J116").Value = 2149.625
$ws.Range("K116").Value = 1409.2059
$ws.Range("L116").Value = 2149.625
$ws.Range("M116").Value = 884.7941000000001
$ws.Range("N116").Value = -6737.625

$ws.Range("H122").Value = 2256.1035
$ws.Range("I122").Value = 1838.8572
$ws.Range("J122").Value = 3351.375
$ws.Range("K122").Value = 5516.571599999999
$ws.Range("L122").Value = 10054.125
$ws.Range("M122").Value = -3066.571599999999
$ws.Range("N122").Value = -14954.125

$ws.Range("H132").Value = 23686.309
$ws.Range("I132").Value = 1388.5
$ws.Range("J132").Value = 127742.75
$ws.Range("K132").Value = 4165.5
$ws.Range("L132").Value = 383228.25
$ws.Range("M132").Value = -1635.5
$ws.Range("N132").Value = -388288.25

$ws.Range("H135").Value = 53514.668
$ws.Range("J135").Value = 53514.668
$ws.Range("L135").Value = 53514.668
$ws.Range("N135").Value = -63654.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1550.238
$ws.Range("I3").Value = 1409.2059
$ws.Range("J3").Value = 2149.625
$ws.Range("K3").Value = 1409.2059
$ws.Range("L3").Value = 2149.625
$ws.Range("M3").Value = -1295.2059
$ws.Range("N3").Value = -2377.625

$ws.Range("H22").Value = 555855.4399999999
$ws.Range("I22").Value = 714485.7
$ws.Range("J22").Value = 649.5
$ws.Range("K22").Value = 714485.7
$ws.Range("L22").Value = 649.5
$ws.Range("M22").Value = -714312.7
$ws.Range("N22").Value = -995.5

$ws.Range("H86").Value = 1450.2122
$ws.Range("I86").Value = 1340.25
$ws.Range("J86").Value = 1619.3846
$ws.Range("K86").Value = 1340.25
$ws.Range("L86").Value = 1619.3846
$ws.Range("M86").Value = -217.25
$ws.Range("N86").Value = -3865.3846

$ws.Range("H89").Value = 1450.2122
$ws.Range("I89").Value = 1340.25
$ws.Range("J89").Value = 1619.3846
$ws.Range("K89").Value = 6701.25
$ws.Range("L89").Value = 8096.923000000001
$ws.Range("M89").Value = -1085.25
$ws.Range("N89").Value = -19328.923

$ws.Range("H105").Value = 2519.7144
$ws.Range("I105").Value = 2440.9333
$ws.Range("K105").Value = 2440.9333
$ws.Range("M105").Value = -693.9333000000001

$ws.Range("H107").Value = 823.90625
$ws.Range("J107").Value = 1001.2143
$ws.Range("L107").Value = 1001.2143
$ws.Range("N107").Value = -4841.2143

$ws.Range("H134").Value = 25280.09
$ws.Range("I134").Value = 26859.62
$ws.Range("K134").Value = 80578.86
$ws.Range("M134").Value = -78043.86

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 253.6
$ws.Range("I7").Value = 267
$ws.Range("K7").Value = 267
$ws.Range("M7").Value = -154

$ws.Range("H58").Value = 17141.258
$ws.Range("I58").Value = 1127.2916
$ws.Range("J58").Value = 72046.28999999999
$ws.Range("K58").Value = 1127.2916
$ws.Range("L58").Value = 72046.28999999999
$ws.Range("M58").Value = -924.2916
$ws.Range("N58").Value = -72452.28999999999

$ws.Range("H74").Value = 27904.924
$ws.Range("J74").Value = 27904.924
$ws.Range("L74").Value = 27904.924
$ws.Range("N74").Value = -29652.924

$ws.Range("H77").Value = 27904.924
$ws.Range("J77").Value = 27904.924
$ws.Range("L77").Value = 83714.772
$ws.Range("N77").Value = -92450.772

$ws.Range("H122").Value = 1306.7858
$ws.Range("J122").Value = 1469.6154
$ws.Range("L122").Value = 4408.8462
$ws.Range("N122").Value = -9308.8462

$ws.Range("H132").Value = 12289.6455
$ws.Range("I132").Value = 14508.105
$ws.Range("J132").Value = 3859.5
$ws.Range("K132").Value = 43524.315
$ws.Range("L132").Value = 11578.5
$ws.Range("M132").Value = -40994.315
$ws.Range("N132").Value = -16638.5

$ws.Range("H134").Value = 902.03705
$ws.Range("I134").Value = 767.75
$ws.Range("J134").Value = 1285.7142
$ws.Range("K134").Value = 2303.25
$ws.Range("L134").Value = 3857.1426
$ws.Range("M134").Value = 231.75
$ws.Range("N134").Value = -8927.142599999999

$ws.Range("H136").Value = 17141.258
$ws.Range("I136").Value = 1127.2916
$ws.Range("J136").Value = 72046.28999999999
$ws.Range("K136").Value = 3381.8748
$ws.Range("L136").Value = 216138.87
$ws.Range("M136").Value = -831.8748000000001
$ws.Range("N136").Value = -221238.87

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 8214.429
$ws.Range("J80").Value = 10036.637
$ws.Range("L80").Value = 30109.911
$ws.Range("N80").Value = -31981.911

$ws.Range("H83").Value = 8214.429
$ws.Range("J83").Value = 10036.637
$ws.Range("L83").Value = 90329.73300000001
$ws.Range("N83").Value = -99689.73300000001

$ws.Range("H105").Value = 4000
$ws.Range("J105").Value = 4000
$ws.Range("L105").Value = 12000
$ws.Range("N105").Value = -17242

$ws.Range("H131").Value = 781.8
$ws.Range("J131").Value = 794.6391599999999
$ws.Range("L131").Value = 2383.91748
$ws.Range("N131").Value = -12463.91748

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6000
$ws.Range("I70").Value = 6000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = -5730
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 6000
$ws.Range("I73").Value = 6000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = -5064
$ws.Range("M73").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 40063
$ws.Range("J133").Value = 40063
$ws.Range("L133").Value = 40063
$ws.Range("N133").Value = -45123

$ws.Range("H136").Value = 24294.637
$ws.Range("I136").Value = 32410.812
$ws.Range("K136").Value = 97232.436
$ws.Range("M136").Value = -94682.436

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1299.3334
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 1449
$ws.Range("K96").Value = 1000
$ws.Range("L96").Value = 1449
$ws.Range("M96").Value = 373
$ws.Range("N96").Value = -4195

$ws.Range("H132").Value = 871.7692
$ws.Range("I132").Value = 661.1389
$ws.Range("K132").Value = 1983.4167
$ws.Range("M132").Value = 546.5832999999998
